$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '30.677.38'
Set-TextValue 'E2' '  +1.15%  '

# Row 3
Set-TextValue 'D3' '2.101.59'
Set-TextValue 'E3' '  +5.12%  '

# Row 4
Set-TextValue 'D4' '1.004'
Set-TextValue 'E4' '  +0.27%  '

# Row 5
Set-TextValue 'D5' '330.04'
Set-TextValue 'E5' '  +1.98%  '

# Row 6
Set-TextValue 'E6' '  +0.22%  '

# Row 7
Set-TextValue 'D7' '0.5276'
Set-TextValue 'E7' '  +3.45%  '

# Row 8
Set-TextValue 'D8' '0.4335'
Set-TextValue 'E8' '  +5.07%  '

# Row 9
Set-TextValue 'D9' '0.08910'
Set-TextValue 'E9' '  +2.28%  '

# Row 10
Set-TextValue 'D10' '46.73'
Set-TextValue 'E10' '  +9.50%  '

# Row 11
Set-TextValue 'E11' '  +2.63%  '

# Row 12
Set-TextValue 'D12' '24.56'
Set-TextValue 'E12' '  -0.63%  '

# Row 13
Set-TextValue 'D13' '2.101.65'
Set-TextValue 'E13' '  +5.26%  '

# Row 14
Set-TextValue 'D14' '6.687'
Set-TextValue 'E14' '  +2.26%  '

# Row 15
Set-TextValue 'D15' '7.760'
Set-TextValue 'E15' '  +4.52%  '

# Row 16
Set-TextValue 'D16' '97.20'
Set-TextValue 'E16' '  +3.38%  '

# Row 17
Set-TextValue 'E17' '  +0.19%  '

# Row 18
Set-TextValue 'D18' '0.00001126'
Set-TextValue 'E18' '  +0.99%  '

# Row 19
Set-TextValue 'D19' '0.06657'
Set-TextValue 'E19' '  +2.27%  '

# Row 20
Set-TextValue 'D20' '18.92'
Set-TextValue 'E20' '  +0.04%  '

# Row 21
Set-TextValue 'E21' '  +0.17%  '

# Row 22
Set-TextValue 'E22' '  +1.64%  '

# Row 23
Set-TextValue 'D23' '30.749.59'
Set-TextValue 'E23' '  +1.23%  '

# Row 24
Set-TextValue 'D24' '12.24'
Set-TextValue 'E24' '  +3.56%  '

# Row 25
Set-TextValue 'D25' '2.358.81'
Set-TextValue 'E25' '  +5.37%  '

# Row 26
Set-TextValue 'D26' '2.288'
Set-TextValue 'E26' '  +3.37%  '

# Row 27
Set-TextValue 'D27' '22.45'
Set-TextValue 'E27' '  -0.01%  '

# Row 28
Set-TextValue 'D28' '2.550'
Set-TextValue 'E28' '  +5.70%  '

# Row 29
Set-TextValue 'E29' '  -0.91%  '

# Row 30
Set-TextValue 'D30' '132.53'
Set-TextValue 'E30' '  +0.72%  '

# Row 31
Set-TextValue 'D31' '1.193'
Set-TextValue 'E31' '  +4.51%  '

# Row 32
Set-TextValue 'E32' '  +2.24%  '

# Row 33
Set-TextValue 'D33' '6.128'
Set-TextValue 'E33' '  +1.24%  '

# Row 34
Set-TextValue 'D34' '1.537'
Set-TextValue 'E34' '  +15.20%  '

# Row 35
Set-TextValue 'D35' '3.848'
Set-TextValue 'E35' '  +0.43%  '

# Row 36
Set-TextValue 'D36' '0.02583'
Set-TextValue 'E36' '  +2.77%  '

# Row 37
Set-TextValue 'D37' '9.601'
Set-TextValue 'E37' '  +6.61%  '

# Row 38
Set-TextValue 'D38' '5.507'
Set-TextValue 'E38' '  +2.58%  '

# Row 39
Set-TextValue 'D39' '0.06713'
Set-TextValue 'E39' '  +1.65%  '

# Row 40
Set-TextValue 'B40' 'Algorand'
Set-TextValue 'C40' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D40' '0.2263'
Set-TextValue 'E40' '  +2.85%  '

# Row 41
Set-TextValue 'B41' 'Aptos'
Set-TextValue 'C41' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D41' '12.54'
Set-TextValue 'E41' '  +2.56%  '

# Row 42
Set-TextValue 'D42' '0.6773'
Set-TextValue 'E42' '  +2.45%  '

# Row 43
Set-TextValue 'E43' '  +1.32%  '

# Row 44
Set-TextValue 'E44' '  +0.14%  '

# Row 45
Set-TextValue 'D45' '0.6374'
Set-TextValue 'E45' '  +3.60%  '

# Row 46
Set-TextValue 'D46' '13.97'
Set-TextValue 'E46' '  +2.08%  '

# Row 47
Set-TextValue 'D47' '2.207'
Set-TextValue 'E47' '  +0.18%  '

# Row 48
Set-TextValue 'D48' '3.631'
Set-TextValue 'E48' '  -0.91%  '

# Row 49
Set-TextValue 'E49' '  -0.87%  '

# Row 50
Set-TextValue 'D50' '82.49'
Set-TextValue 'E50' '  +2.75%  '

# Row 51
Set-TextValue 'D51' '1.194'
Set-TextValue 'E51' '  +6.30%  '
